# Fix closing dates / rates that were wrong for algo
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JIND")

# Row 7
$ws.Range("F7").Value = 697.15
$ws.Range("G7").Value = 710
$ws.Range("H7").Value = 691.1
$ws.Range("I7").Value = 706.05
$ws.Range("J7").Value = 697.55

# Row 9
$ws.Range("G9").Value = 700
$ws.Range("H9").Value = 685.05
$ws.Range("I9").Value = 693.85

# Row 10
$ws.Range("G10").Value = 700.85
$ws.Range("H10").Value = 691.1
$ws.Range("I10").Value = 698.35

# Row 11
$ws.Range("G11").Value = 702.6
$ws.Range("H11").Value = 697.05
$ws.Range("I11").Value = 699.85

# Row 12
$ws.Range("G12").Value = 700.55
$ws.Range("H12").Value = 694.8
$ws.Range("I12").Value = 699.3

# Row 13
$ws.Range("G13").Value = 701.5
$ws.Range("H13").Value = 697.4
$ws.Range("I13").Value = 697.4

# Row 14
$ws.Range("G14").Value = 701.2
$ws.Range("H14").Value = 697.6
$ws.Range("I14").Value = 700.7

# Row 15
$ws.Range("G15").Value = 703.45
$ws.Range("H15").Value = 700.4
$ws.Range("I15").Value = 702.45

# Row 16
$ws.Range("G16").Value = 706.7
$ws.Range("H16").Value = 702.45
$ws.Range("I16").Value = 706.25

# Row 17
$ws.Range("G17").Value = 710
$ws.Range("H17").Value = 704.45
$ws.Range("I17").Value = 706.5

# Row 18
$ws.Range("G18").Value = 709.95
$ws.Range("H18").Value = 703.45
$ws.Range("I18").Value = 705.9

# Row 19
$ws.Range("G19").Value = 708
$ws.Range("H19").Value = 702.5
$ws.Range("I19").Value = 707.85

# Row 20
$ws.Range("G20").Value = 710
$ws.Range("H20").Value = 704.65
$ws.Range("I20").Value = 705.4

# Row 21
$ws.Range("G21").Value = 707.15
$ws.Range("H21").Value = 704.5
$ws.Range("I21").Value = 707
